$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP addresses in columns A (host) and C (ip) for rows 2 and 3
$ws.Range("A2").Value = "10.10.10.1"
$ws.Range("C2").Value = "10.10.10.1"
$ws.Range("A3").Value = "10.10.10.2"
$ws.Range("C3").Value = "10.10.10.2"

# Update the current selection to match the saved state
$ws.Range("B4").Select()
